$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Feuil1")

# Update the "Part" column for both "Button" rows (front & rear brake systems)
# to the more descriptive "Rotor button".
$ws.Range("C6").Value = "Rotor button"
$ws.Range("C12").Value = "Rotor button"

# Move the active selection to C13, matching the author's cursor position
# at the time of the edit.
$ws.Range("C13").Select()
